# Append-refresh of the scraped Lancers job listing sheet.
# A new scrape ran at 2025-10-14 01:17:23 JST: the 9 previously-seen postings are kept
# (with their "取得日時" timestamp refreshed) and 4 newly-seen postings are merged in,
# with the full list re-sorted by "優先度スコア" (column G) descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (タイトル) widened from 42 to 51 characters.
# NOTE: the ColumnWidth COM property reads ~0.83 narrower than the character width stored
# in the saved XML "width" attribute, so we compensate to land exactly on 51.
$ws.Columns.Item(2).ColumnWidth = 50.17

# Drop all existing hyperlinks up front. Rows are about to be rewritten in place (several
# shift down to make room for the new postings), and re-adding them fresh below guarantees
# every hyperlink ends up attached to the correct, final cell.
$ws.Hyperlinks.Delete()

# Row 2: 【注目】ChatGPTを活用した薬歴アプリ開発の依頼
$ws.Cells.Item(2,1).Value = '2025-10-14 01:17:23'
$ws.Cells.Item(2,2).Value = '【注目】ChatGPTを活用した薬歴アプリ開発の依頼'
$ws.Cells.Item(2,3).Value = 'システム開発'
$ws.Cells.Item(2,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(2,5).Value = '期限情報なし'
$ws.Cells.Item(2,6).Value = 'https://www.lancers.jp/work/detail/5412417'
$ws.Hyperlinks.Add($ws.Cells.Item(2,6), 'https://www.lancers.jp/work/detail/5412417')
$ws.Cells.Item(2,6).Style = "Hyperlink"
$ws.Cells.Item(2,7).Value = 398
$ws.Cells.Item(2,8).Value = '🔥GPT,ChatGPT ◆開発 ◇アプリ'

# Row 3: 【GAS開発】配送状況管理の自動化を依頼します
$ws.Cells.Item(3,1).Value = '2025-10-14 01:17:23'
$ws.Cells.Item(3,2).Value = '【GAS開発】配送状況管理の自動化を依頼します'
$ws.Cells.Item(3,3).Value = 'システム開発'
$ws.Cells.Item(3,4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(3,5).Value = '期限情報なし'
$ws.Cells.Item(3,6).Value = 'https://www.lancers.jp/work/detail/5412306'
$ws.Hyperlinks.Add($ws.Cells.Item(3,6), 'https://www.lancers.jp/work/detail/5412306')
$ws.Cells.Item(3,6).Style = "Hyperlink"
$ws.Cells.Item(3,7).Value = 170
$ws.Cells.Item(3,8).Value = '◆開発,自動化 ◇管理'

# Row 4: Amazon商品を自動抽出してBASEに出品するツール開発(スクレイピング機能)
$ws.Cells.Item(4,1).Value = '2025-10-14 01:17:23'
$ws.Cells.Item(4,2).Value = 'Amazon商品を自動抽出してBASEに出品するツール開発(スクレイピング機能)'
$ws.Cells.Item(4,3).Value = 'システム開発'
$ws.Cells.Item(4,4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(4,5).Value = '期限情報なし'
$ws.Cells.Item(4,6).Value = 'https://www.lancers.jp/work/detail/5412467'
$ws.Hyperlinks.Add($ws.Cells.Item(4,6), 'https://www.lancers.jp/work/detail/5412467')
$ws.Cells.Item(4,6).Style = "Hyperlink"
$ws.Cells.Item(4,7).Value = 168
$ws.Cells.Item(4,8).Value = '◆ツール,開発'

# Row 5: 海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動
$ws.Cells.Item(5,1).Value = '2025-10-14 01:17:23'
$ws.Cells.Item(5,2).Value = '海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)'
$ws.Cells.Item(5,3).Value = 'システム開発'
$ws.Cells.Item(5,4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(5,5).Value = '期限情報なし'
$ws.Cells.Item(5,6).Value = 'https://www.lancers.jp/work/detail/5251319'
$ws.Hyperlinks.Add($ws.Cells.Item(5,6), 'https://www.lancers.jp/work/detail/5251319')
$ws.Cells.Item(5,6).Style = "Hyperlink"
$ws.Cells.Item(5,7).Value = 135
$ws.Cells.Item(5,8).Value = '◆ツール,スクレイピング ◇サイト'

# Row 6: 【急募】クリニック向け内視鏡画像システム開発の依頼
$ws.Cells.Item(6,1).Value = '2025-10-14 01:17:23'
$ws.Cells.Item(6,2).Value = '【急募】クリニック向け内視鏡画像システム開発の依頼'
$ws.Cells.Item(6,3).Value = 'システム開発'
$ws.Cells.Item(6,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(6,5).Value = '期限情報なし'
$ws.Cells.Item(6,6).Value = 'https://www.lancers.jp/work/detail/5412233'
$ws.Hyperlinks.Add($ws.Cells.Item(6,6), 'https://www.lancers.jp/work/detail/5412233')
$ws.Cells.Item(6,6).Style = "Hyperlink"
$ws.Cells.Item(6,7).Value = 125
$ws.Cells.Item(6,8).Value = '◆開発,システム開発'

# Row 7: 【急募】onedrive上のexcelで自動化システム構築依頼
$ws.Cells.Item(7,1).Value = '2025-10-14 01:17:23'
$ws.Cells.Item(7,2).Value = '【急募】onedrive上のexcelで自動化システム構築依頼'
$ws.Cells.Item(7,3).Value = 'システム開発'
$ws.Cells.Item(7,4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(7,5).Value = '期限情報なし'
$ws.Cells.Item(7,6).Value = 'https://www.lancers.jp/work/detail/5412194'
$ws.Hyperlinks.Add($ws.Cells.Item(7,6), 'https://www.lancers.jp/work/detail/5412194')
$ws.Cells.Item(7,6).Style = "Hyperlink"
$ws.Cells.Item(7,7).Value = 95
$ws.Cells.Item(7,8).Value = '◆自動化'

# Row 8: 【フルリモート】WordPress開発スタッフ募集
$ws.Cells.Item(8,1).Value = '2025-10-14 01:17:23'
$ws.Cells.Item(8,2).Value = '【フルリモート】WordPress開発スタッフ募集'
$ws.Cells.Item(8,3).Value = 'システム開発'
$ws.Cells.Item(8,4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(8,5).Value = '期限情報なし'
$ws.Cells.Item(8,6).Value = 'https://www.lancers.jp/work/detail/5407811'
$ws.Hyperlinks.Add($ws.Cells.Item(8,6), 'https://www.lancers.jp/work/detail/5407811')
$ws.Cells.Item(8,6).Style = "Hyperlink"
$ws.Cells.Item(8,7).Value = 88
$ws.Cells.Item(8,8).Value = '◆開発 ○WordPress'

# Row 9: 【急募】スタートアップ向けプロダクト開発のパートナー募集
$ws.Cells.Item(9,1).Value = '2025-10-14 01:17:23'
$ws.Cells.Item(9,2).Value = '【急募】スタートアップ向けプロダクト開発のパートナー募集'
$ws.Cells.Item(9,3).Value = 'システム開発'
$ws.Cells.Item(9,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(9,5).Value = '期限情報なし'
$ws.Cells.Item(9,6).Value = 'https://www.lancers.jp/work/detail/5412179'
$ws.Hyperlinks.Add($ws.Cells.Item(9,6), 'https://www.lancers.jp/work/detail/5412179')
$ws.Cells.Item(9,6).Style = "Hyperlink"
$ws.Cells.Item(9,7).Value = 75
$ws.Cells.Item(9,8).Value = '◆開発'

# Row 10: FileMaker開発
$ws.Cells.Item(10,1).Value = '2025-10-14 01:17:23'
$ws.Cells.Item(10,2).Value = 'FileMaker開発'
$ws.Cells.Item(10,3).Value = 'システム開発'
$ws.Cells.Item(10,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(10,5).Value = '期限情報なし'
$ws.Cells.Item(10,6).Value = 'https://www.lancers.jp/work/detail/5412487'
$ws.Hyperlinks.Add($ws.Cells.Item(10,6), 'https://www.lancers.jp/work/detail/5412487')
$ws.Cells.Item(10,6).Style = "Hyperlink"
$ws.Cells.Item(10,7).Value = 68
$ws.Cells.Item(10,8).Value = '◆開発'

# Row 11: 初回 ポケパラの自動いいね等の開発
$ws.Cells.Item(11,1).Value = '2025-10-14 01:17:23'
$ws.Cells.Item(11,2).Value = '初回 ポケパラの自動いいね等の開発'
$ws.Cells.Item(11,3).Value = 'システム開発'
$ws.Cells.Item(11,4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(11,5).Value = '期限情報なし'
$ws.Cells.Item(11,6).Value = 'https://www.lancers.jp/work/detail/5412453'
$ws.Hyperlinks.Add($ws.Cells.Item(11,6), 'https://www.lancers.jp/work/detail/5412453')
$ws.Cells.Item(11,6).Style = "Hyperlink"
$ws.Cells.Item(11,7).Value = 63
$ws.Cells.Item(11,8).Value = '◆開発'

# Row 12: 【音声コマンド起動】超小型・低電力レコーダーのプロトタイプ開発
$ws.Cells.Item(12,1).Value = '2025-10-14 01:17:23'
$ws.Cells.Item(12,2).Value = '【音声コマンド起動】超小型・低電力レコーダーのプロトタイプ開発'
$ws.Cells.Item(12,3).Value = 'システム開発'
$ws.Cells.Item(12,4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(12,5).Value = '期限情報なし'
$ws.Cells.Item(12,6).Value = 'https://www.lancers.jp/work/detail/5412261'
$ws.Hyperlinks.Add($ws.Cells.Item(12,6), 'https://www.lancers.jp/work/detail/5412261')
$ws.Cells.Item(12,6).Style = "Hyperlink"
$ws.Cells.Item(12,7).Value = 60
$ws.Cells.Item(12,8).Value = '◆開発'

# Row 13: 微生物の特定と分類を行いたく、画像解析の専門家を探しています!(急いでません!)
$ws.Cells.Item(13,1).Value = '2025-10-14 01:17:23'
$ws.Cells.Item(13,2).Value = '微生物の特定と分類を行いたく、画像解析の専門家を探しています!(急いでません!)'
$ws.Cells.Item(13,3).Value = 'システム開発'
$ws.Cells.Item(13,4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(13,5).Value = '期限情報なし'
$ws.Cells.Item(13,6).Value = 'https://www.lancers.jp/work/detail/5411887'
$ws.Hyperlinks.Add($ws.Cells.Item(13,6), 'https://www.lancers.jp/work/detail/5411887')
$ws.Cells.Item(13,6).Style = "Hyperlink"
$ws.Cells.Item(13,7).Value = 18
$ws.Cells.Item(13,8).ClearContents()

# Row 14: LINE公式(Lステップ)のリッチメニューの構築
$ws.Cells.Item(14,1).Value = '2025-10-14 01:17:23'
$ws.Cells.Item(14,2).Value = 'LINE公式(Lステップ)のリッチメニューの構築'
$ws.Cells.Item(14,3).Value = 'システム開発'
$ws.Cells.Item(14,4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(14,5).Value = '期限情報なし'
$ws.Cells.Item(14,6).Value = 'https://www.lancers.jp/work/detail/5412357'
$ws.Hyperlinks.Add($ws.Cells.Item(14,6), 'https://www.lancers.jp/work/detail/5412357')
$ws.Cells.Item(14,6).Style = "Hyperlink"
$ws.Cells.Item(14,7).Value = 10
$ws.Cells.Item(14,8).ClearContents()
